$wb = $excel.ActiveWorkbook

# Copy the "Croatia" sheet (template layout) to the end of the workbook to
# create the new "Greece" market sheet.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market-specific values on the new sheet.
$greece.Range("B4").Value = "NGC-4119/T3168"
$greece.Range("B2").Value = "Greece Market"

# Croatia is no longer the active tab; its selection becomes the whole sheet.
[void]$croatia.Cells.Select()

# Make Greece the active/selected sheet, matching the authored selection.
[void]$greece.Activate()
[void]$greece.Range("H20").Select()
